{"js": "const pairs = [\n  [\"2025-12-10 Wednesday\", \"2025-12-11 Thursday\"],\n  [\"314\u00f75=62, 4\", \"543\u00f75=108, 3\"],\n  [\"784\u00f79=87, 1\", \"434\u00f76=72, 2\"],\n  [\"957\u00f73=319, 0\", \"122\u00f76=20, 2\"],\n  [\"803\u00f76=133, 5\", \"997\u00f77=142, 3\"],\n  [\"507\u00f79=56, 3\", \"913\u00f77=130, 3\"],\n  [\"759\u00f78=94, 7\", \"275\u00f77=39, 2\"],\n  [\"571\u00f76=95, 1\", \"720\u00f75=144, 0\"],\n  [\"557\u00f77=79, 4\", \"468\u00f72=234, 0\"],\n  [\"495\u00f75=99, 0\", \"849\u00f78=106, 1\"],\n  [\"158\u00f73=52, 2\", \"503\u00f74=125, 3\"],\n  [\"933\u00f77=133, 2\", \"525\u00f72=262, 1\"],\n  [\"879\u00f78=109, 7\", \"363\u00f73=121, 0\"],\n  [\"527\u00f77=75, 2\", \"123\u00f72=61, 1\"],\n  [\"224\u00f79=24, 8\", \"993\u00f75=198, 3\"],\n  [\"269\u00f79=29, 8\", \"712\u00f72=356, 0\"],\n  [\"821\u00f79=91, 2\", \"317\u00f77=45, 2\"],\n  [\"174\u00f72=87, 0\", \"367\u00f79=40, 7\"],\n  [\"481\u00f73=160, 1\", \"533\u00f75=106, 3\"],\n  [\"833\u00f72=416, 1\", \"654\u00f77=93, 3\"],\n  [\"589\u00f74=147, 1\", \"899\u00f77=128, 3\"],\n  [\"973\u00f76=162, 1\", \"252\u00f76=42, 0\"],\n  [\"736\u00f73=245, 1\", \"217\u00f78=27, 1\"],\n  [\"369\u00f77=52, 5\", \"347\u00f78=43, 3\"],\n  [\"724\u00f73=241, 1\", \"752\u00f76=125, 2\"],\n  [\"780\u00f78=97, 4\", \"247\u00f77=35, 2\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Not found: \" + oldText);\n  }\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-Text($doc, $oldText, $newText) {\n  $find = $doc.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $find.Forward = $true\n  $find.Wrap = 1  # wdFindContinue\n  $find.Format = $false\n  $find.MatchCase = $true\n  $find.MatchWholeWord = $false\n  $find.MatchWildcards = $false\n  $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n\nReplace-Text $d \"2025-12-10 Wednesday\" \"2025-12-11 Thursday\"\nReplace-Text $d \"314\u00f75=62, 4\" \"543\u00f75=108, 3\"\nReplace-Text $d \"784\u00f79=87, 1\" \"434\u00f76=72, 2\"\nReplace-Text $d \"957\u00f73=319, 0\" \"122\u00f76=20, 2\"\nReplace-Text $d \"803\u00f76=133, 5\" \"997\u00f77=142, 3\"\nReplace-Text $d \"507\u00f79=56, 3\" \"913\u00f77=130, 3\"\nReplace-Text $d \"759\u00f78=94, 7\" \"275\u00f77=39, 2\"\nReplace-Text $d \"571\u00f76=95, 1\" \"720\u00f75=144, 0\"\nReplace-Text $d \"557\u00f77=79, 4\" \"468\u00f72=234, 0\"\nReplace-Text $d \"495\u00f75=99, 0\" \"849\u00f78=106, 1\"\nReplace-Text $d \"158\u00f73=52, 2\" \"503\u00f74=125, 3\"\nReplace-Text $d \"933\u00f77=133, 2\" \"525\u00f72=262, 1\"\nReplace-Text $d \"879\u00f78=109, 7\" \"363\u00f73=121, 0\"\nReplace-Text $d \"527\u00f77=75, 2\" \"123\u00f72=61, 1\"\nReplace-Text $d \"224\u00f79=24, 8\" \"993\u00f75=198, 3\"\nReplace-Text $d \"269\u00f79=29, 8\" \"712\u00f72=356, 0\"\nReplace-Text $d \"821\u00f79=91, 2\" \"317\u00f77=45, 2\"\nReplace-Text $d \"174\u00f72=87, 0\" \"367\u00f79=40, 7\"\nReplace-Text $d \"481\u00f73=160, 1\" \"533\u00f75=106, 3\"\nReplace-Text $d \"833\u00f72=416, 1\" \"654\u00f77=93, 3\"\nReplace-Text $d \"589\u00f74=147, 1\" \"899\u00f77=128, 3\"\nReplace-Text $d \"973\u00f76=162, 1\" \"252\u00f76=42, 0\"\nReplace-Text $d \"736\u00f73=245, 1\" \"217\u00f78=27, 1\"\nReplace-Text $d \"369\u00f77=52, 5\" \"347\u00f78=43, 3\"\nReplace-Text $d \"724\u00f73=241, 1\" \"752\u00f76=125, 2\"\nReplace-Text $d \"780\u00f78=97, 4\" \"247\u00f77=35, 2\"\n"}
